# Add the new "Tyden 34" worksheet at the end of the workbook and populate
# it with the weekly roster + shift-time data.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Týden 34"

# Shift start time (07:00) stored as a day-fraction, formatted HH:MM.
$ws.Range("D7").Value = 0.2916666666666667
$ws.Range("D7").NumberFormat = "HH:MM"

# Employee roster with a per-row hour total (formatted with 2 decimals).
$names = @(
    @(9,  "Čáp Jakub"),
    @(10, "Horčička Jiří"),
    @(11, "Hromý Erik"),
    @(12, "Kužel Andrej"),
    @(13, "Mlynář Roman"),
    @(14, "Nastoupil Ladislav"),
    @(15, "Winkler Jan"),
    @(16, "Štrauf Daniel")
)

foreach ($entry in $names) {
    $row = $entry[0]
    $name = $entry[1]
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 5).Value = 10
    $ws.Cells.Item($row, 5).NumberFormat = "0.00"
}

# Week-ending date, first registered as yyyy-mm-dd then switched to
# DD.MM.YYYY (mirrors the source workbook, which keeps the now-unused
# yyyy-mm-dd numFmt entry around from that earlier step).
$ws.Range("D80").Value = 45888
$ws.Range("D80").NumberFormat = "yyyy-mm-dd"
$ws.Range("D80").NumberFormat = "DD.MM.YYYY"
